# feature: enviar certificado por email
# Adds an "Email" column next to the existing "Nome" column on the first
# sheet ("Abertura Versando a Cidade"), drops the stray leading "Data"
# header row, and fills in the known email addresses for the first three
# listed people.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the current (single) column out of the way to column B, preserving
# its width/bestFit metadata, and free up column A for the "Nome" list.
$ws.Columns.Item(1).Insert()

$black = 0
$fontName = "Arial"
$fontSize = 11

function Set-HeaderCell($addr, $text) {
    $c = $ws.Range($addr)
    $c.Value = $text
    $c.Font.Name = $fontName
    $c.Font.Size = $fontSize
    $c.Font.Bold = $true
    $c.Font.Color = $black
    $c.HorizontalAlignment = -4108   # xlCenter
}

function Set-DataCell($addr, $text) {
    $c = $ws.Range($addr)
    $c.Value = $text
    $c.Font.Name = $fontName
    $c.Font.Size = $fontSize
    $c.Font.Bold = $false
    $c.Font.Color = $black
    $c.HorizontalAlignment = -4131   # xlLeft
}

function Set-QuotePrefixDataCell($addr) {
    # Same look as the "Pessoa N" / email data cells, but left empty with a
    # leading apostrophe (quote-prefixed empty text).
    $c = $ws.Range($addr)
    $c.Font.Name = $fontName
    $c.Font.Size = $fontSize
    $c.Font.Bold = $false
    $c.Font.Color = $black
    $c.HorizontalAlignment = -4131   # xlLeft
    $c.Value = "'"                  # lone apostrophe -> stored as empty quoted text
}

function Set-QuotePrefixGeneralCell($addr) {
    # Same look as the plain blank rows further down (default font/general
    # alignment), but quote-prefixed empty text instead of a true blank.
    $c = $ws.Range($addr)
    $c.HorizontalAlignment = 1       # xlGeneral
    $c.Value = "'"                  # lone apostrophe -> stored as empty quoted text
}

function Set-BlankGeneralCell($addr) {
    $c = $ws.Range($addr)
    $c.HorizontalAlignment = 1       # xlGeneral
}

function Set-BlankLeftCell($addr) {
    $c = $ws.Range($addr)
    $c.Font.Name = $fontName
    $c.Font.Size = $fontSize
    $c.Font.Bold = $false
    $c.Font.Color = $black
    $c.HorizontalAlignment = -4131   # xlLeft
}

# Row 1: headers
Set-HeaderCell "A1" "Nome"
Set-HeaderCell "B1" "Email"

# Rows 2-4: known people + their emails
Set-DataCell "A2" "Pessoa 1"
Set-DataCell "B2" "allanfernds@gmail.com"

Set-DataCell "A3" "Pessoa 2"
Set-DataCell "B3" "allanweik@gmail.com"

Set-DataCell "A4" "Pessoa 3"
Set-DataCell "B4" "alanfernandes.mm@gmail.com"

# Rows 5-6: placeholder rows entered with a leading apostrophe (empty text)
Set-QuotePrefixGeneralCell "A5"
Set-QuotePrefixDataCell "B5"
Set-QuotePrefixGeneralCell "A6"
Set-QuotePrefixDataCell "B6"

# Rows 7-36: remaining blank rows to keep the same usable range size
# (row 36 is brand new, so it needs the same row height as the rest)
$ws.Rows.Item(36).RowHeight = 19.5
for ($r = 7; $r -le 36; $r++) {
    Set-BlankGeneralCell ("A" + $r)
    Set-BlankLeftCell ("B" + $r)
}

# Column widths: "Nome" column narrower, "Email" column keeps the width
# inherited from the original single column (already preserved by Insert()).
$ws.Columns.Item(1).ColumnWidth = 31.2

Write-Output "done"
